$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Difficulty column (G) values for existing questions to simple
# whole-number difficulty ratings 1-9.
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 4
$ws.Range("G6").Value = 5
$ws.Range("G7").Value = 6
$ws.Range("G8").Value = 7
$ws.Range("G9").Value = 8
$ws.Range("G10").Value = 9

# Add a new quiz question row (row 11) - the "change difficulty" joke question.
$ws.Range("A11").Value = "Is this the hardest question in the quiz?  9.99"
$ws.Range("B11").Value = "Yes"
$ws.Range("C11").Value = "N"
$ws.Range("D11").Value = "No"
$ws.Range("E11").Value = "No"
$ws.Range("F11").Value = "Yes"
$ws.Range("G11").Value = 9.99

# Match the saved selection state from the authored workbook.
[void]$ws.Range("G9").Select()
